# Updated symbol list on Fri Jan  6 14:56:00 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest scraped coinranking.com values. All cells in these columns are
# stored as literal text (e.g. "256.58", "0.25%") rather than numbers, so
# each value is written with a leading apostrophe to force Excel to keep it
# as text instead of auto-converting it to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'256.58"
$ws.Range("E2").Value  = "'0.25%"

$ws.Range("D3").Value  = "'26.55"
$ws.Range("E3").Value  = "'-1.48%"

$ws.Range("D4").Value  = "'4.643"
$ws.Range("E4").Value  = "'0.11%"

$ws.Range("D5").Value  = "'0.05921"
$ws.Range("E5").Value  = "'0.63%"

$ws.Range("D6").Value  = "'6.600"
$ws.Range("E6").Value  = "'-0.63%"

$ws.Range("D7").Value  = "'0.8559"
$ws.Range("E7").Value  = "'-1.43%"

$ws.Range("D8").Value  = "'0.9139"
$ws.Range("E8").Value  = "'-2.36%"

$ws.Range("D9").Value  = "'0.1377"
$ws.Range("E9").Value  = "'-1.70%"

$ws.Range("D10").Value = "'0.04450"
$ws.Range("E10").Value = "'16.70%"

$ws.Range("D11").Value = "'0.06999"
$ws.Range("E11").Value = "'-1.18%"

$ws.Range("D12").Value = "'0.03030"
$ws.Range("E12").Value = "'-5.61%"

$ws.Range("D13").Value = "'0.09102"
$ws.Range("E13").Value = "'-1.46%"

$ws.Range("D14").Value = "'0.001527"
$ws.Range("E14").Value = "'-1.10%"

$ws.Range("D15").Value = "'0.0006030"
$ws.Range("E15").Value = "'0.48%"

$ws.Range("D16").Value = "'0.006060"
$ws.Range("E16").Value = "'0.90%"

$ws.Range("D17").Value = "'3.464"
$ws.Range("E17").Value = "'-1.47%"

$ws.Range("D18").Value = "'3.130"
$ws.Range("E18").Value = "'-1.97%"

$ws.Range("E19").Value = "'-2.37%"

$ws.Range("D20").Value = "'0.3078"
$ws.Range("E20").Value = "'0.15%"

$ws.Range("D21").Value = "'0.1287"
$ws.Range("E21").Value = "'0.39%"

$ws.Range("D22").Value = "'3.894"
$ws.Range("E22").Value = "'1.19%"

$ws.Range("D23").Value = "'0.04209"
$ws.Range("E23").Value = "'-0.27%"

$ws.Range("E24").Value = "'-0.13%"

$ws.Range("D25").Value = "'0.004611"
$ws.Range("E25").Value = "'8.01%"

$ws.Range("E26").Value = "'0.01%"

$ws.Range("D27").Value = "'0.0001715"
$ws.Range("E27").Value = "'13.75%"

$ws.Range("D40").Value = "'0.03797"
$ws.Range("E40").Value = "'-0.44%"

$ws.Range("D41").Value = "'0.006271"
$ws.Range("E41").Value = "'59.85%"

$ws.Range("D42").Value = "'0.1097"
$ws.Range("E42").Value = "'-0.14%"

$ws.Range("D43").Value = "'0.002200"
$ws.Range("E43").Value = "'-3.81%"

$ws.Range("E44").Value = "'24.56%"

$ws.Range("D45").Value = "'0.00005108"
$ws.Range("E45").Value = "'-6.58%"

$ws.Range("E46").Value = "'-0.01%"

$ws.Range("D47").Value = "'0.05001"
$ws.Range("E47").Value = "'-16.95%"

$ws.Range("E48").Value = "'10,469.03%"

$ws.Range("E49").Value = "'-0.01%"

$ws.Range("E50").Value = "'-0.01%"
